$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 37.5
$ws.Range("I5").Value = 37.5
$ws.Range("K5").Value = 37.5
$ws.Range("M5").Value = 77.5
# Row 28
$ws.Range("H28").Value = 468
$ws.Range("I28").Value = 468
$ws.Range("K28").Value = 468
$ws.Range("M28").Value = 17
# Row 70
$ws.Range("H70").Value = 1019.0909
$ws.Range("J70").Value = 1298.75
$ws.Range("L70").Value = 3896.25
$ws.Range("N70").Value = -4436.25
# Row 73
$ws.Range("H73").Value = 1019.0909
$ws.Range("J73").Value = 1298.75
$ws.Range("L73").Value = 3896.25
$ws.Range("N73").Value = -5768.25
# Row 100
$ws.Range("H100").Value = 66521.48
$ws.Range("I100").Value = 72823.86
$ws.Range("J100").Value = 56717.777
$ws.Range("K100").Value = 72823.86
$ws.Range("L100").Value = 56717.777
$ws.Range("M100").Value = -72282.86
$ws.Range("N100").Value = -57799.777
# Row 112
$ws.Range("H112").Value = 2513.6667
$ws.Range("I112").Value = 5500
$ws.Range("J112").Value = 1916.4
$ws.Range("K112").Value = 16500
$ws.Range("L112").Value = 5749.200000000001
$ws.Range("M112").Value = -15392
$ws.Range("N112").Value = -7965.200000000001
# Row 141
$ws.Range("H141").Value = 8781.440000000001
$ws.Range("I141").Value = 9196.799999999999
$ws.Range("K141").Value = 27590.4
$ws.Range("M141").Value = -22410.4

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1979.2
$ws.Range("I74").Value = 678.625
$ws.Range("K74").Value = 678.625
$ws.Range("M74").Value = 195.375
# Row 77
$ws.Range("H77").Value = 1979.2
$ws.Range("I77").Value = 678.625
$ws.Range("K77").Value = 3393.125
$ws.Range("M77").Value = 974.875
# Row 132
$ws.Range("H132").Value = 1880.2
$ws.Range("I132").Value = 1205.8096
$ws.Range("K132").Value = 3617.4288
$ws.Range("M132").Value = -1087.4288
# Row 135
$ws.Range("H135").Value = 90164.164
$ws.Range("J135").Value = 90164.164
$ws.Range("L135").Value = 90164.164
$ws.Range("N135").Value = -100304.164

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 4140.5
$ws.Range("I36").Value = 4188
$ws.Range("J36").Value = 3998
$ws.Range("K36").Value = 4188
$ws.Range("L36").Value = 3998
$ws.Range("M36").Value = -3654
$ws.Range("N36").Value = -5066
# Row 86
$ws.Range("H86").Value = 8245.117
$ws.Range("I86").Value = 10937.2
$ws.Range("K86").Value = 10937.2
$ws.Range("M86").Value = -9814.200000000001
# Row 89
$ws.Range("H89").Value = 8245.117
$ws.Range("I89").Value = 10937.2
$ws.Range("K89").Value = 54686
$ws.Range("M89").Value = -49070
# Row 105
$ws.Range("H105").Value = 64458.055
$ws.Range("I105").Value = 111897.6
$ws.Range("J105").Value = 5158.625
$ws.Range("K105").Value = 111897.6
$ws.Range("L105").Value = 5158.625
$ws.Range("M105").Value = -110150.6
$ws.Range("N105").Value = -8652.625
# Row 107
$ws.Range("H107").Value = 2655.375
$ws.Range("I107").Value = 2989
$ws.Range("K107").Value = 2989
$ws.Range("M107").Value = -1069

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 767283.5
$ws.Range("I6").Value = 767283.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 767283.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -767170.5
$ws.Range("N6").ClearContents()
# Row 7
$ws.Range("H7").Value = 267.86667
$ws.Range("I7").Value = 219.8
$ws.Range("K7").Value = 219.8
$ws.Range("M7").Value = -106.8
# Row 51
$ws.Range("H51").Value = 7833.3335
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 59
$ws.Range("H59").Value = 39778.5
$ws.Range("J59").Value = 39778.5
$ws.Range("L59").Value = 39778.5
$ws.Range("N59").Value = -42068.5
# Row 60
$ws.Range("H60").Value = 5044.875
# Row 61
$ws.Range("H61").Value = 7833.3335
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 68
$ws.Range("H68").Value = 39687.6
$ws.Range("J68").Value = 40541.777
$ws.Range("L68").Value = 40541.777
$ws.Range("N68").Value = -42039.777
# Row 71
$ws.Range("H71").Value = 39687.6
$ws.Range("J71").Value = 40541.777
$ws.Range("L71").Value = 121625.331
$ws.Range("N71").Value = -129113.331
# Row 74
$ws.Range("H74").Value = 60233.668
$ws.Range("J74").Value = 65534.6
$ws.Range("L74").Value = 65534.6
$ws.Range("N74").Value = -67282.60000000001
# Row 77
$ws.Range("H77").Value = 60233.668
$ws.Range("J77").Value = 65534.6
$ws.Range("L77").Value = 196603.8
$ws.Range("N77").Value = -205339.8
# Row 107
$ws.Range("H107").Value = 8625.467000000001
$ws.Range("I107").Value = 12208.3
$ws.Range("K107").Value = 12208.3
$ws.Range("M107").Value = -10288.3
# Row 122
$ws.Range("H122").Value = 9227.5625
$ws.Range("I122").Value = 13428.3
$ws.Range("K122").Value = 40284.89999999999
$ws.Range("M122").Value = -37834.89999999999
# Row 132
$ws.Range("H132").Value = 1734.7188
$ws.Range("I132").Value = 1517.0667
$ws.Range("K132").Value = 4551.2001
$ws.Range("M132").Value = -2021.2001
# Row 141
$ws.Range("H141").Value = 182791.12
$ws.Range("J141").Value = 198652.19
$ws.Range("L141").Value = 198652.19
$ws.Range("N141").Value = -209012.19

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 11183.667
$ws.Range("J68").Value = 13240.5
$ws.Range("L68").Value = 39721.5
$ws.Range("N68").Value = -41343.5
# Row 71
$ws.Range("H71").Value = 11183.667
$ws.Range("J71").Value = 13240.5
$ws.Range("L71").Value = 119164.5
$ws.Range("N71").Value = -127276.5
# Row 129
$ws.Range("H129").Value = 23810064
$ws.Range("I129").Value = 489.375
$ws.Range("J129").Value = 55556164
$ws.Range("K129").Value = 1468.125
$ws.Range("L129").Value = 166668492
$ws.Range("M129").Value = 3531.875
$ws.Range("N129").Value = -166678492

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 13337.733
$ws.Range("I80").Value = 15465.917
$ws.Range("J80").Value = 4825
$ws.Range("K80").Value = 15465.917
$ws.Range("L80").Value = 4825
$ws.Range("M80").Value = -14467.917
$ws.Range("N80").Value = -6821
# Row 83
$ws.Range("H83").Value = 13337.733
$ws.Range("I83").Value = 15465.917
$ws.Range("J83").Value = 4825
$ws.Range("K83").Value = 77329.58499999999
$ws.Range("L83").Value = 24125
$ws.Range("M83").Value = -72337.58499999999
$ws.Range("N83").Value = -34109
# Row 122
$ws.Range("H122").Value = 10987.556
$ws.Range("I122").Value = 7576.857
$ws.Range("K122").Value = 22730.571
$ws.Range("M122").Value = -20280.571
# Row 132
$ws.Range("H132").Value = 2427.7812
$ws.Range("I132").Value = 2325.7778
$ws.Range("K132").Value = 6977.3334
$ws.Range("M132").Value = -4447.3334

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1523.8649
$ws.Range("I16").Value = 1338.8
$ws.Range("K16").Value = 1338.8
$ws.Range("M16").Value = -1168.8
# Row 40
$ws.Range("H40").Value = 20746.967
$ws.Range("I40").Value = 26088.15
$ws.Range("K40").Value = 26088.15
$ws.Range("M40").Value = -25952.15
# Row 46
$ws.Range("H46").Value = 2507.1365
$ws.Range("I46").Value = 754.4286
$ws.Range("K46").Value = 754.4286
$ws.Range("M46").Value = -566.4286
# Row 136
$ws.Range("H136").Value = 4222.7437
$ws.Range("I136").Value = 2979.8572
$ws.Range("J136").Value = 5672.778
$ws.Range("K136").Value = 8939.571599999999
$ws.Range("L136").Value = 17018.334
$ws.Range("M136").Value = -6389.571599999999
$ws.Range("N136").Value = -22118.334

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 6356.7896
$ws.Range("I81").Value = 7286.125
$ws.Range("J81").Value = 1400.3334
$ws.Range("K81").Value = 14572.25
$ws.Range("L81").Value = 2800.6668
$ws.Range("M81").Value = -13511.25
$ws.Range("N81").Value = -4922.6668
# Row 84
$ws.Range("H84").Value = 6356.7896
$ws.Range("I84").Value = 7286.125
$ws.Range("J84").Value = 1400.3334
$ws.Range("K84").Value = 72861.25
$ws.Range("L84").Value = 14003.334
$ws.Range("M84").Value = -67557.25
$ws.Range("N84").Value = -24611.334
# Row 132
$ws.Range("H132").Value = 8353.214
$ws.Range("I132").Value = 9177.761
$ws.Range("K132").Value = 27533.283
$ws.Range("M132").Value = -25003.283
# Row 136
$ws.Range("H136").Value = 297650.25
$ws.Range("I136").Value = 309316.28
$ws.Range("K136").Value = 927948.8400000001
$ws.Range("M136").Value = -925398.8400000001
